# Update countries & provincias Spain
#
# 1) Swap displayed country names for rows 60/61 (Suiza <-> Uzbekistan)
# 2) Refresh the "Datos actualizados ..." timestamp in A1
# 3) Refresh the COVID counters (Casos totales / Nuevos casos / Casos activos /
#    Recuperados / Casos criticos / Muertes hoy / Muertes) for the countries
#    whose stats moved: Ucrania, Suiza, Uzbekistan, Kirguistan, El Salvador,
#    Australia, Tailandia.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Suiza / Uzbekistan swap (rows 60 and 61, column A) -----------------
$ws.Range("A60").Value = "Uzbekistan"
$ws.Range("A61").Value = "Suiza"

# --- 2) Timestamp update (row 1) --------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 07:55"

# --- 3) Counter refresh ------------------------------------------------------
# Row 27 - Ucrania
$ws.Range("B27").Value = 151859
$ws.Range("C27").Value = 3103
$ws.Range("D27").Value = 68346
$ws.Range("E27").Value = 80365
$ws.Range("G27").Value = 72
$ws.Range("H27").Value = 3148

# Row 60 - Uzbekistan (after the swap above)
$ws.Range("B60").Value = 46498
$ws.Range("C60").Value = 338
$ws.Range("D60").Value = 43023
$ws.Range("E60").Value = 3093
$ws.Range("G60").Value = 5
$ws.Range("H60").Value = 382

# Row 61 - Suiza (after the swap above)
$ws.Range("B61").Value = 46239
$ws.Range("D61").Value = 38100
$ws.Range("E61").Value = 6119
$ws.Range("H61").Value = 2020

# Row 64 - Kirguistan
$ws.Range("B64").Value = 44828
$ws.Range("C64").Value = 67
$ws.Range("D64").Value = 40779
$ws.Range("E64").Value = 2986

# Row 74 - El Salvador
$ws.Range("D74").Value = 17446
$ws.Range("E74").Value = 8545
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 782

# Row 75 - Australia
$ws.Range("B75").Value = 26607
$ws.Range("C75").Value = 42
$ws.Range("D75").Value = 23329
$ws.Range("E75").Value = 2475

# Row 128 - Tailandia
$ws.Range("B128").Value = 3466
$ws.Range("C128").Value = 5
$ws.Range("E128").Value = 96
